# Update cryptocurrency Price (D) and Volume(1h) (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$dCell = $ws.Range("D2")
$dCell.NumberFormat = "@"
$dCell.Value = '64.121.78'
$dCell.Style = "Normal"
$ws.Range("E2").Value = '  +0.92%  '

# Row 3
$dCell = $ws.Range("D3")
$dCell.NumberFormat = "@"
$dCell.Value = '3.075.94'
$dCell.Style = "Normal"
$ws.Range("E3").Value = '  +0.18%  '

# Row 4
$ws.Range("E4").Value = '  +0.05%  '

# Row 5
$dCell = $ws.Range("D5")
$dCell.NumberFormat = "@"
$dCell.Value = '557.43'
$dCell.Style = "Normal"
$ws.Range("E5").Value = '  +1.10%  '

# Row 6
$dCell = $ws.Range("D6")
$dCell.NumberFormat = "@"
$dCell.Value = '144.09'
$dCell.Style = "Normal"
$ws.Range("E6").Value = '  +0.98%  '

# Row 7
$ws.Range("E7").Value = '  +0.08%  '

# Row 8
$dCell = $ws.Range("D8")
$dCell.NumberFormat = "@"
$dCell.Value = '3.073.81'
$dCell.Style = "Normal"
$ws.Range("E8").Value = '  +0.31%  '

# Row 9
$dCell = $ws.Range("D9")
$dCell.NumberFormat = "@"
$dCell.Value = '0.507'
$dCell.Style = "Normal"
$ws.Range("E9").Value = '  +0.72%  '

# Row 10
$dCell = $ws.Range("D10")
$dCell.NumberFormat = "@"
$dCell.Value = '0.155'
$dCell.Style = "Normal"
$ws.Range("E10").Value = '  +2.23%  '

# Row 11
$dCell = $ws.Range("D11")
$dCell.NumberFormat = "@"
$dCell.Value = '6.07'
$dCell.Style = "Normal"
$ws.Range("E11").Value = '  -7.15%  '

# Row 12
$dCell = $ws.Range("D12")
$dCell.NumberFormat = "@"
$dCell.Value = '0.472'
$dCell.Style = "Normal"
$ws.Range("E12").Value = '  +3.10%  '

# Row 13
$dCell = $ws.Range("D13")
$dCell.NumberFormat = "@"
$dCell.Value = '0.0000229'
$dCell.Style = "Normal"
$ws.Range("E13").Value = '  +0.36%  '

# Row 14
$dCell = $ws.Range("D14")
$dCell.NumberFormat = "@"
$dCell.Value = '34.99'
$dCell.Style = "Normal"
$ws.Range("E14").Value = '  +0.08%  '

# Row 15
$dCell = $ws.Range("D15")
$dCell.NumberFormat = "@"
$dCell.Value = '3.589.97'
$dCell.Style = "Normal"
$ws.Range("E15").Value = '  +0.66%  '

# Row 16
$dCell = $ws.Range("D16")
$dCell.NumberFormat = "@"
$dCell.Value = '64.143.32'
$dCell.Style = "Normal"
$ws.Range("E16").Value = '  +1.04%  '

# Row 17
$dCell = $ws.Range("D17")
$dCell.NumberFormat = "@"
$dCell.Value = '3.079.64'
$dCell.Style = "Normal"
$ws.Range("E17").Value = '  +0.22%  '

# Row 18
$ws.Range("E18").Value = '  +1.13%  '

# Row 19
$ws.Range("E19").Value = '  -1.00%  '

# Row 20
$dCell = $ws.Range("D20")
$dCell.NumberFormat = "@"
$dCell.Value = '479.57'
$dCell.Style = "Normal"
$ws.Range("E20").Value = '  -1.22%  '

# Row 21
$dCell = $ws.Range("D21")
$dCell.NumberFormat = "@"
$dCell.Value = '14.01'
$dCell.Style = "Normal"
$ws.Range("E21").Value = '  +0.93%  '

# Row 22
$ws.Range("E22").Value = '  -0.03%  '

# Row 23
$dCell = $ws.Range("D23")
$dCell.NumberFormat = "@"
$dCell.Value = '7.52'
$dCell.Style = "Normal"
$ws.Range("E23").Value = '  +2.90%  '

# Row 24
$dCell = $ws.Range("D24")
$dCell.NumberFormat = "@"
$dCell.Value = '14.09'
$dCell.Style = "Normal"
$ws.Range("E24").Value = '  +10.11%  '

# Row 25
$dCell = $ws.Range("D25")
$dCell.NumberFormat = "@"
$dCell.Value = '81.33'
$dCell.Style = "Normal"
$ws.Range("E25").Value = '  +0.25%  '

# Row 26
$ws.Range("E26").Value = '  +0.05%  '

# Row 27
$ws.Range("E27").Value = '  +0.17%  '

# Row 28
$dCell = $ws.Range("D28")
$dCell.NumberFormat = "@"
$dCell.Value = '7.97'
$dCell.Style = "Normal"
$ws.Range("E28").Value = '  +0.92%  '

# Row 29
$ws.Range("E29").Value = '  +1.47%  '

# Row 30
$dCell = $ws.Range("D30")
$dCell.NumberFormat = "@"
$dCell.Value = '1.00'
$dCell.Style = "Normal"
$ws.Range("E30").Value = '  +0.14%  '

# Row 31
$dCell = $ws.Range("D31")
$dCell.NumberFormat = "@"
$dCell.Value = '26.25'
$dCell.Style = "Normal"
$ws.Range("E31").Value = '  +0.01%  '

# Row 32
$ws.Range("E32").Value = '  -2.19%  '

# Row 33
$dCell = $ws.Range("D33")
$dCell.NumberFormat = "@"
$dCell.Value = '2.46'
$dCell.Style = "Normal"
$ws.Range("E33").Value = '  +0.28%  '

# Row 34
$dCell = $ws.Range("D34")
$dCell.NumberFormat = "@"
$dCell.Value = '5.62'
$dCell.Style = "Normal"
$ws.Range("E34").Value = '  -1.59%  '

# Row 35
$dCell = $ws.Range("D35")
$dCell.NumberFormat = "@"
$dCell.Value = '6.19'
$dCell.Style = "Normal"
$ws.Range("E35").Value = '  +2.95%  '

# Row 36
$dCell = $ws.Range("D36")
$dCell.NumberFormat = "@"
$dCell.Value = '55.31'
$dCell.Style = "Normal"
$ws.Range("E36").Value = '  -0.48%  '

# Row 37
$ws.Range("E37").Value = '  +2.07%  '

# Row 38
$dCell = $ws.Range("D38")
$dCell.NumberFormat = "@"
$dCell.Value = '2.93'
$dCell.Style = "Normal"
$ws.Range("E38").Value = '  +13.83%  '

# Row 39
$dCell = $ws.Range("D39")
$dCell.NumberFormat = "@"
$dCell.Value = '437.90'
$dCell.Style = "Normal"
$ws.Range("E39").Value = '  -6.37%  '

# Row 40
$dCell = $ws.Range("D40")
$dCell.NumberFormat = "@"
$dCell.Value = '0.0808'
$dCell.Style = "Normal"
$ws.Range("E40").Value = '  -2.07%  '

# Row 41
$dCell = $ws.Range("D41")
$dCell.NumberFormat = "@"
$dCell.Value = '2.957.66'
$dCell.Style = "Normal"
$ws.Range("E41").Value = '  -2.75%  '

# Row 42
$dCell = $ws.Range("D42")
$dCell.NumberFormat = "@"
$dCell.Value = '8.19'
$dCell.Style = "Normal"
$ws.Range("E42").Value = '  -0.70%  '

# Row 43
$ws.Range("E43").Value = '  -4.73%  '

# Row 44
$dCell = $ws.Range("D44")
$dCell.NumberFormat = "@"
$dCell.Value = '28.11'
$dCell.Style = "Normal"
$ws.Range("E44").Value = '  +0.89%  '

# Row 45
$ws.Range("E45").Value = '  +1.39%  '

# Row 46
$ws.Range("E46").Value = '  +0.03%  '

# Row 47
$dCell = $ws.Range("D47")
$dCell.NumberFormat = "@"
$dCell.Value = '2.12'
$dCell.Style = "Normal"
$ws.Range("E47").Value = '  +3.64%  '

# Row 48
$ws.Range("E48").Value = '  +1.05%  '

# Row 49
$ws.Range("E49").Value = '  +0.80%  '

# Row 50
$dCell = $ws.Range("D50")
$dCell.NumberFormat = "@"
$dCell.Value = '117.48'
$dCell.Style = "Normal"
$ws.Range("E50").Value = '  +0.72%  '

# Row 51
$ws.Range("E51").Value = '  -0.65%  '
